# "Version 2." -> "Version 1."
#
# Character layout of the paragraph before editing:
#   0-5  "Versi"   (run 1)
#   5-7  "on"      (run 2, w:rsidR="008807B2")
#   7-9  " 2"      (run 3, w:rsidR="008807B2")
#   9-10 "."       (run 4)
#   10   paragraph mark
#
# Edits are applied back-to-front so earlier character offsets stay valid.

$d = $word.ActiveDocument

# 1) Drop the trailing "." run entirely.
$d.Range(9, 10).Text = ""

# 2) Turn " 2" into " 1." (run survives, text content changes).
$d.Range(7, 9).Text = " 1."

# 3) Merge the split "Versi" + "on" runs into a single "Version" run:
#    first delete "on" ...
$d.Range(5, 7).Text = ""
#    ... then rewrite "Versi" as "Version" so Word collapses to one run.
$d.Range(0, 5).Text = "Version"
